# "fix the new dungeon info"
# - Clears the stray Quest value in F23 (密林迷宫 / forestmaze) that no longer applies.
# - Adds three new dungeon sub-rows (13020021-13020023) describing the new
#   village dungeon (登拉克峡谷 / 13010004): 村落入口, 村中心, 议事厅.
# - Grows the "表1" structured table and the related ranges (dimension,
#   conditional formatting) to cover the 3 new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Row 23 (密林迷宫): the Quest column no longer has a value.
# ---------------------------------------------------------------------------
$ws.Range("F23").ClearContents()

# ---------------------------------------------------------------------------
# 2. Grow the structured table ("表1") by three rows (A3:T27 -> A3:T30).
# ---------------------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$lo.ListRows.Add() | Out-Null
$lo.ListRows.Add() | Out-Null
$lo.ListRows.Add() | Out-Null

# Copy the formatting of the last existing data row onto the new rows so the
# new cells pick up the same styles (s="14"/"9"/"12" ...) as the rest of the
# table.
$ws.Range("A27:T27").Copy()
$ws.Range("A28:T30").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3. Fill in the three new dungeon rows.
# ---------------------------------------------------------------------------
# Row 28: 村落入口 (village entrance)
$ws.Cells.Item(28, 1).Value = 13020021
$ws.Cells.Item(28, 2).Value = "村落入口"
$ws.Cells.Item(28, 3).Value = 3
$ws.Cells.Item(28, 4).Value = 5
$ws.Cells.Item(28, 5).Value = 13010004
$ws.Cells.Item(28, 8).Value = "trees;4"
$ws.Cells.Item(28, 16).Value = "viliage1"
$ws.Cells.Item(28, 17).Value = "viliage1"

# Row 29: 村中心 (village center)
$ws.Cells.Item(29, 1).Value = 13020022
$ws.Cells.Item(29, 2).Value = "村中心"
$ws.Cells.Item(29, 3).Value = 3
$ws.Cells.Item(29, 4).Value = 5
$ws.Cells.Item(29, 5).Value = 13010004
$ws.Cells.Item(29, 16).Value = "viliage2"
$ws.Cells.Item(29, 17).Value = "viliage2"

# Row 30: 议事厅 (council hall)
$ws.Cells.Item(30, 1).Value = 13020023
$ws.Cells.Item(30, 2).Value = "议事厅"
$ws.Cells.Item(30, 3).Value = 3
$ws.Cells.Item(30, 4).Value = 5
$ws.Cells.Item(30, 5).Value = 13010004
$ws.Cells.Item(30, 16).Value = "viliage3"
$ws.Cells.Item(30, 17).Value = "viliage3"

# ---------------------------------------------------------------------------
# 4. Also re-save row 23's H column (QuestDungeon) and Q23 (TilePath) with the
#    same value they already held, so the workbook stores a fresh reference
#    for them (text content is unchanged).
# ---------------------------------------------------------------------------
$ws.Cells.Item(23, 8).Value = "trees;4|manflower;2|portal;1|sandland;2|cliff;2|colorpool;1|barn;1"
$ws.Cells.Item(23, 17).Value = "forestmaze"

# ---------------------------------------------------------------------------
# 5. Extend the conditional formatting range that covered I4:N27 so that it
#    now covers the new rows too (I4:N30). The other rule (J4 only) is left
#    untouched.
# ---------------------------------------------------------------------------
$bigRule = $ws.Range("I4:N27").FormatConditions.Item(1)
$bigRule.ModifyAppliesToRange($ws.Range("I4:N30"))

# ---------------------------------------------------------------------------
# 6. Update the selection to match the edited cell.
# ---------------------------------------------------------------------------
$ws.Range("F23").Select()
